$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row at row 17 for "TRITONE 200MG 30 TAB" (alphabetically
#    it sits between "SUGARLO PLUS ..." (row 16) and "VIDROP ..." (old row 17,
#    which becomes row 18 after the insert).
# ---------------------------------------------------------------------------
$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value = 11

$c = $ws.Range("C17"); $c.NumberFormat = "@"; $c.Value = "TRITONE 200MG 30 TAB"
$c = $ws.Range("H17"); $c.NumberFormat = "@"; $c.Value = "1:2"
$c = $ws.Range("L17"); $c.NumberFormat = "@"; $c.Value = "1"
$c = $ws.Range("N17"); $c.NumberFormat = "@"; $c.Value = "123.00"
$c = $ws.Range("P17"); $c.NumberFormat = "@"; $c.Value = "-40.5900"
$c = $ws.Range("Q17"); $c.NumberFormat = "@"; $c.Value = "0:-1"

# Re-apply the look of the surrounding table (borders/fonts/number formats)
# from the row immediately below, which still carries the original styling.
$ws.Range("A18:Q18").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4122)

$ws.Rows.Item(17).RowHeight = 25.5
$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

# ---------------------------------------------------------------------------
# 2) Insert a second new data row at row 20 for "ZYMAGALLIN 30 TAB" (sits
#    between "WATER FOR INJECTION ..." (row 19) and the Arabic entries that
#    start with "بلاستر مترسيلك ..." (row 20, shifting to row 21)).
# ---------------------------------------------------------------------------
$ws.Rows.Item(20).Insert()

$ws.Range("A20").Value = 14

$c = $ws.Range("C20"); $c.NumberFormat = "@"; $c.Value = "ZYMAGALLIN 30 TAB"
$c = $ws.Range("H20"); $c.NumberFormat = "@"; $c.Value = "2:1"
$c = $ws.Range("L20"); $c.NumberFormat = "@"; $c.Value = "1"
$c = $ws.Range("N20"); $c.NumberFormat = "@"; $c.Value = "48.00"
$c = $ws.Range("P20"); $c.NumberFormat = "@"; $c.Value = "-15.8400"
$c = $ws.Range("Q20"); $c.NumberFormat = "@"; $c.Value = "0:-1"

$ws.Range("A21:Q21").Copy()
$ws.Range("A20:Q20").PasteSpecial(-4122)

$ws.Rows.Item(20).RowHeight = 24.75
$ws.Range("A20:B20").Merge()
$ws.Range("C20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()
$ws.Range("N20:O20").Merge()

# ---------------------------------------------------------------------------
# 3) Renumber the trailing rows (A column index) so the sequence stays
#    1..17 across the now-17 data rows.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = 15
$ws.Range("A22").Value = 16
$ws.Range("A23").Value = 17

# ---------------------------------------------------------------------------
# 4) Update the running total (sum of column P) which is a literal value,
#    not a formula, and now lives on row 24 after the two inserts above.
# ---------------------------------------------------------------------------
$ws.Range("P24").Value = 702.585

# ---------------------------------------------------------------------------
# 5) Refresh the "printed at" timestamp in the footer (now row 25).
# ---------------------------------------------------------------------------
$ws.Range("K25").Value = "Saturday, 27 September, 2025 11:52 AM"
